$d = $word.ActiveDocument

# --- "Viena limena detalizetaks apraksts" heading: Heading 2 -> Heading 3 ---
$d.Paragraphs(11).Style = "Heading 3"

# --- Remove the stray run containing a single space right before the
#     "NK miss2.jpg" picture. ---
$shp = $d.InlineShapes(1)
$picRange = $shp.Range
$spaceStart = $picRange.Start - 1
$spaceRng = $d.Range($spaceStart, $picRange.Start)
if ($spaceRng.Text -eq " ") {
    $spaceRng.Delete()
}

# --- Drop the two trailing "materials" sub-bullets under "Stalazu
#     buvesana" -> "Ja vajadzigs iezime citu materialu" (the fixed-length /
#     drag-and-drop notes). ---
$start = $d.Paragraphs(23).Range.Start
$end = $d.Paragraphs(24).Range.End
$d.Range($start, $end).Delete()

# --- Drop the whole "Kopejam materialu ... / Materialiem ... / Materialus
#     spelitajs unlock-o ... / Koka delis ... / Terauda balsts ... / Virve
#     - ???" block of scrap/brainstorm bullets. ---
$start = $d.Paragraphs(16).Range.Start
$end = $d.Paragraphs(21).Range.End
$d.Range($start, $end).Delete()

# --- Remove the leading stray space before "animacija, kur speletajam
#     pieskir ..." ---
$p = $d.Paragraphs(24)
$leadStart = $p.Range.Start
$leadRng = $d.Range($leadStart, $leadStart + 1)
if ($leadRng.Text -eq " ") {
    $leadRng.Delete()
}

Write-Output ("Paragraphs: " + $d.Paragraphs.Count)
